$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "sides_theroost"

# --- Expand the table by one row (was A1:G5 -> A1:G6) ---
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# --- Row 1 (header) ---
$ws.Range("A1").Value = "ItemName"
$ws.Range("B1").Value = "Ingredients"
$ws.Range("C1").Value = "Allergens"
$ws.Range("D1").Value = "LocalIngredients"
$ws.Range("E1").Value = "Diet"
$ws.Range("F1").Value = "Nutrition Label"
$ws.Range("G1").Value = "LeaveEmpty"

# --- Row 2 : Fries ---
$ws.Range("A2").Value = "Fries"
$ws.Range("B2").Value = "McCain Gold Crisp Fries[Potatoes, canola oil, wheat flour, modified corn starch, corn flour, salt, autolyzed yeast, baking powder, sodium phosphate, dextrose, modified cellulose, colour]"
$ws.Range("C2").Value = "Wheat."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "VGN,DF"
$ws.Range("F2").Value = "Fries"

# --- Row 3 : Cheese Curds ---
$ws.Range("A3").Value = "Cheese Curds"
$ws.Range("B3").Value = "29% MF Milk / Enzymes / Salt"
$ws.Range("C3").Value = "Milk."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "GF,VEG"
$ws.Range("F3").Value = "placeholder"

# --- Row 4 : Chef's Salad ---
$ws.Range("A4").Value = "Chef's Salad"
$ws.Range("B4").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("C4").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "placeholder"
# B4/C4 used to carry the wrap-text style (old row3/row4 content) - drop it, row 4 is plain now
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Style = "Normal"

# --- Row 5 : Coleslaw ---
$ws.Range("A5").Value = "Coleslaw"
$ws.Range("B5").Value = "Green Cabbage / Red Cabbage / Shredded Carrot / Coleslaw Dressing"
$ws.Range("C5").Value = "Eggs, mustard."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "GF,VEG"
$ws.Range("F5").Value = "House_Salad"

# --- Row 6 : Daily Soup ---
$ws.Range("A6").Value = "Daily Soup"
$ws.Range("B6").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("C6").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "placeholder"

# --- Styling ---
# New centered/bold-ish nutrition icon style on A2
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108

# Wrap-text style (existing style index 1) re-applied on the ingredient/allergen cells
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true

# --- Row heights ---
# Rows 2, 3 & 4 lose their old custom heights -> back to standard
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
# Row 5 keeps the old 15.75pt custom height that used to live on row 4
$ws.Rows.Item(5).RowHeight = 15.75
